$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 9.352362000000001
$ws.Range("H2").Value = 28.057086
$ws.Range("I2").Value = 0.3814309586590714
$ws.Range("J2").Value = 0.3814309586590714
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.899038333333333
$ws.Range("N2").Value = 5.697115
$ws.Range("Q2").Value = 17.76049394521
$ws.Range("R2").Value = 159.84444550689
$ws.Range("S2").Value = 0.3814309586590714
$ws.Range("T2").Value = 0.3814309586590714

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 8.765748333333333
$ws.Range("H3").Value = 26.297245
$ws.Range("I3").Value = 0.3575062417544884
$ws.Range("J3").Value = 0.3575062417544883
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.899038333333333
$ws.Range("N3").Value = 5.697115
$ws.Range("Q3").Value = 16.64649210535278
$ws.Range("R3").Value = 149.818428948175
$ws.Range("S3").Value = 0.3575062417544884
$ws.Range("T3").Value = 0.3575062417544883

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.331995666666666
$ws.Range("H4").Value = 9.995987
$ws.Range("I4").Value = 0.1358936171829681
$ws.Range("J4").Value = 0.1358936171829681
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.899038333333333
$ws.Range("N4").Value = 5.697115
$ws.Range("Q4").Value = 6.327587497500555
$ws.Range("R4").Value = 56.948287477505
$ws.Range("S4").Value = 0.1358936171829681
$ws.Range("T4").Value = 0.1358936171829681

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.069041666666667
$ws.Range("H5").Value = 9.207125000000001
$ws.Range("I5").Value = 0.1251691824034721
$ws.Range("J5").Value = 0.1251691824034721
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.899038333333333
$ws.Range("N5").Value = 5.697115
$ws.Range("Q5").Value = 5.828227771597223
$ws.Range("R5").Value = 52.45404994437501
$ws.Range("S5").Value = 0.1251691824034721
$ws.Range("T5").Value = 0.1251691824034721
